# texts.xlsx edit — "Found right version of program"
#
# 1) Typography sheet: remove the stray "AI_value / SourceSansPro-Regular.ttf"
#    typography row (row 8, columns B:J) that wasn't a real typography entry.
#    The K:O "Wildcard Ranges" table living in the same physical row is left
#    untouched.
# 2) Translation sheet: the "AI_value" text id had incorrectly been reused as
#    a Typography Name value for several translation rows — fix those rows
#    back to "Typography_01", and correct the Text Id / Alignment / GB values
#    for the whole block (rows 5-27) to the right program version's content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Typography sheet — drop the bogus row (B8:J8)
# ---------------------------------------------------------------------------
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Range("B8:J8").ClearContents()

# ---------------------------------------------------------------------------
# 2) Translation sheet — rewrite rows 5-27 (Text Id / Typography Name /
#    Alignment / GB). Column F (Direction) is unchanged throughout.
# ---------------------------------------------------------------------------
$wsTranslation = $wb.Worksheets.Item("Translation")

function Set-TranslationRow($row, $textId, $typographyName, $alignment, $gb) {
    $wsTranslation.Range("B$row").Value2 = $textId
    $wsTranslation.Range("C$row").Value2 = $typographyName
    $wsTranslation.Range("D$row").Value2 = $alignment
    $wsTranslation.Range("E$row").Value2 = $gb
}

Set-TranslationRow 5  'SingleUseId2'  'Typography_01' 'Left'   '<value>; '
Set-TranslationRow 6  'TextId2'       'Typography_00' 'Center' '<number>'
Set-TranslationRow 7  'TextId3'       'Typography_01' 'Center' '<number>'
Set-TranslationRow 8  'SingleUseId15' 'Typography_01' 'Left'   '0'
Set-TranslationRow 9  'TextId4'       'Typography_02' 'Center' 'New Text'
Set-TranslationRow 10 'SingleUseId25' 'Typography_01' 'Left'   '<>'
Set-TranslationRow 11 'SingleUseId27' 'Typography_01' 'Left'   '<value>'
Set-TranslationRow 12 'SingleUseId30' 'Typography_01' 'Left'   '<value>'
Set-TranslationRow 13 'SingleUseId31' 'Typography_01' 'Left'   'admin'
Set-TranslationRow 14 'SingleUseId33' 'Typography_01' 'Left'   'admin'
Set-TranslationRow 15 'SingleUseId35' 'Typography_01' 'Center' '<value>'
Set-TranslationRow 16 'SingleUseId36' 'Typography_01' 'Left'   '000.000'
Set-TranslationRow 17 'SingleUseId53' 'Typography_01' 'Center' '<value>'
Set-TranslationRow 18 'SingleUseId54' 'Typography_01' 'Left'   '000.000'
Set-TranslationRow 19 'SingleUseId55' 'Typography_01' 'Center' '<value>'
Set-TranslationRow 20 'SingleUseId56' 'Typography_01' 'Left'   '000.000'
Set-TranslationRow 21 'SingleUseId57' 'Typography_01' 'Center' '<value>'
Set-TranslationRow 22 'SingleUseId58' 'Typography_01' 'Left'   '000.000'
Set-TranslationRow 23 'SingleUseId70' 'Typography_00' 'Center' '<> %'
Set-TranslationRow 24 'SingleUseId71' 'Typography_00' 'Center' '<> %'
Set-TranslationRow 25 'SingleUseId72' 'Typography_01' 'Left'   'DD/MM/YYYY'
Set-TranslationRow 26 'SingleUseId79' 'Typography_01' 'Left'   '<>'
Set-TranslationRow 27 'SingleUseId80' 'Typography_01' 'Left'   'DD/MM/YYYY'

Write-Output "edit applied"
